$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Header date (A1): 45406 (2024-04-24) -> 45436 (2024-05-24)
$ws.Range("A1").Value = 45436

# Price list entries (D22, D23): 1497.908 -> 2950.798
$ws.Range("D22").Value = 2950.798
$ws.Range("D23").Value = 2950.798
